# Update cryptos list - price & volume(1h) refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    # Force the cell to stay a plain text value (matches the source
    # workbook's inlineStr cells) regardless of whether $value happens to
    # look like a number, and without leaving the cell's style altered.
    $ws.Range($cell).Value = "'" + $value
    $ws.Range($cell).Style = "Normal"
}

Set-TextCell "D2" "29.208.94"
$ws.Range("E2").Value = "  +0.32%  "

Set-TextCell "D3" "1.834.88"
$ws.Range("E3").Value = "  -0.10%  "

$ws.Range("E4").Value = "  +0.02%  "

Set-TextCell "D5" "242.03"
$ws.Range("E5").Value = "  +0.74%  "

Set-TextCell "D6" "0.6662"
$ws.Range("E6").Value = "  -2.28%  "

Set-TextCell "D8" "0.07423"
$ws.Range("E8").Value = "  -0.64%  "

Set-TextCell "D9" "0.2934"
$ws.Range("E9").Value = "  -1.82%  "

Set-TextCell "D10" "22.92"
$ws.Range("E10").Value = "  -1.07%  "

Set-TextCell "D11" "0.07764"
$ws.Range("E11").Value = "  +1.38%  "

Set-TextCell "D12" "1.835.12"
$ws.Range("E12").Value = "  +0.09%  "

Set-TextCell "D13" "5.002"
$ws.Range("E13").Value = "  -0.40%  "

Set-TextCell "D14" "0.6682"
$ws.Range("E14").Value = "  -1.40%  "

$ws.Range("E15").Value = "  -4.19%  "

Set-TextCell "D16" "6.121"
$ws.Range("E16").Value = "  -0.69%  "

Set-TextCell "D17" "0.000008378"
$ws.Range("E17").Value = "  +1.13%  "

Set-TextCell "D18" "29.189.80"
$ws.Range("E18").Value = "  +0.52%  "

Set-TextCell "D19" "2.079.38"
$ws.Range("E19").Value = "  +1.66%  "

Set-TextCell "D20" "228.28"
$ws.Range("E20").Value = "  -0.16%  "

Set-TextCell "D21" "12.47"

$ws.Range("E22").Value = "  +0.22%  "

$ws.Range("E23").Value = "  -2.57%  "

Set-TextCell "D24" "0.9997"
$ws.Range("E24").Value = "  +0.06%  "

Set-TextCell "D25" "159.13"
$ws.Range("E25").Value = "  -1.15%  "

Set-TextCell "D26" "0.1406"
$ws.Range("E26").Value = "  -2.43%  "

Set-TextCell "D27" "8.619"
$ws.Range("E27").Value = "  -1.11%  "

Set-TextCell "D28" "18.00"
$ws.Range("E28").Value = "  -0.19%  "

Set-TextCell "D29" "1.514"
$ws.Range("E29").Value = "  +0.76%  "

Set-TextCell "D30" "4.114"
$ws.Range("E30").Value = "  -3.28%  "

Set-TextCell "D31" "4.046"
$ws.Range("E31").Value = "  -2.08%  "

Set-TextCell "D33" "0.05286"
$ws.Range("E33").Value = "  -2.17%  "

Set-TextCell "D34" "1.864"
$ws.Range("E34").Value = "  +0.26%  "

Set-TextCell "D35" "0.7460"
$ws.Range("E35").Value = "  -0.67%  "

$ws.Range("E36").Value = "  +0.76%  "

Set-TextCell "D37" "2.650"
$ws.Range("E37").Value = "  -1.25%  "

Set-TextCell "D38" "1.306.89"
$ws.Range("E38").Value = "  +0.47%  "

$ws.Range("E39").Value = "  -0.87%  "

Set-TextCell "D40" "2.736"
$ws.Range("E40").Value = "  +0.75%  "

Set-TextCell "D41" "0.9354"
$ws.Range("E41").Value = "  -0.06%  "

Set-TextCell "D42" "5.883"
$ws.Range("E42").Value = "  -2.80%  "

Set-TextCell "D43" "0.08365"
$ws.Range("E43").Value = "  +1.25%  "

$ws.Range("E45").Value = "  -2.46%  "

Set-TextCell "D46" "1.980.63"
$ws.Range("E46").Value = "  +0.59%  "

Set-TextCell "D47" "0.5147"
$ws.Range("E47").Value = "  -0.57%  "

$ws.Range("E48").Value = "  -0.37%  "

$ws.Range("E49").Value = "  -0.96%  "

Set-TextCell "D50" "63.03"
$ws.Range("E50").Value = "  -1.35%  "

$ws.Range("E51").Value = "  -0.91%  "
